# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets,
# to reflect refreshed counts captured at commit 456a3b4 (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1199
$ws1.Range("F4").Value  = 14497
$ws1.Range("F5").Value  = 17233
$ws1.Range("F7").Value  = 144
$ws1.Range("F8").Value  = 52
$ws1.Range("F20").Value = 1311
$ws1.Range("F22").Value = 75
$ws1.Range("F24").Value = 7
$ws1.Range("F25").Value = 7096
$ws1.Range("F26").Value = 977
$ws1.Range("F29").Value = 33
$ws1.Range("F31").Value = 49
$ws1.Range("F32").Value = 5828
$ws1.Range("F33").Value = 131
$ws1.Range("F36").Value = 4997
$ws1.Range("F37").Value = 30

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1199
$ws4.Range("F4").Value  = 14497
$ws4.Range("F5").Value  = 17233
$ws4.Range("F7").Value  = 144
$ws4.Range("F8").Value  = 52
$ws4.Range("F20").Value = 1311
$ws4.Range("F22").Value = 75
$ws4.Range("F26").Value = 7096
$ws4.Range("F27").Value = 977
$ws4.Range("F30").Value = 33
$ws4.Range("F32").Value = 49
$ws4.Range("F34").Value = 5828
$ws4.Range("F35").Value = 131
$ws4.Range("F38").Value = 4997
$ws4.Range("F39").Value = 30
